$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The original export inserted two new columns in front of the old "Is
# Significant" column (F) to make room for "Observed"/"Expected" data, which
# shifts the old F column to H for every data row.
$ws.Range("F2:G2").Insert(-4161)   # xlShiftToRight, cascades to all data rows

# The header row ends up with its labels "out of order" versus the data
# columns below them - F1 reads "Observed", G1 "Expected" and H1 "Is
# Significant" (i.e. the header row never got its old "Is Significant"
# label shifted out of F1 the way every data row below did) - this
# header/data mismatch is the "mess" referenced in the commit message.
$ws.Range("F1").Value = "Observed"
$ws.Range("G1").Value = "Expected"
$ws.Range("H1").Value = "Is Significant"

# Observed / Expected values for each data row (these land in F/G, while the
# old "Is Significant" flag that used to live in F has shifted into H).
$ws.Range("F2").Value = "[ 51 171] ; [16  4]"
$ws.Range("G2").Value = "[ 61.46280992 160.53719008] ; [ 5.53719008 14.46280992]"

$ws.Range("F3").Value = "[ 76 142] ; [19  1]"
$ws.Range("G3").Value = "[ 87.01680672 130.98319328] ; [ 7.98319328 12.01680672]"

$ws.Range("F4").Value = "[ 65 146] ; [16  4]"
$ws.Range("G4").Value = "[ 73.98701299 137.01298701] ; [ 7.01298701 12.98701299]"

$ws.Range("F5").Value = "[  4 220] ; [ 4 16]"
$ws.Range("G5").Value = "[  7.3442623 216.6557377] ; [ 0.6557377 19.3442623]"

$ws.Range("F6").Value = "[ 72 150] ; [16  4]"
$ws.Range("G6").Value = "[ 80.72727273 141.27272727] ; [ 7.27272727 12.72727273]"

$ws.Range("F7").Value = "[102 123] ; [18  2]"
$ws.Range("G7").Value = "[110.20408163 114.79591837] ; [ 9.79591837 10.20408163]"

$ws.Range("F8").Value = "[ 10 214] ; [ 5 14]"
$ws.Range("G8").Value = "[ 13.82716049 210.17283951] ; [ 1.17283951 17.82716049]"

$ws.Range("F9").Value = "[103 119] ; [18  2]"
$ws.Range("G9").Value = "[111. 111.] ; [10. 10.]"

$ws.Range("F10").Value = "[ 81 141] ; [14  3]"
$ws.Range("G10").Value = "[ 88.24267782 133.75732218] ; [ 6.75732218 10.24267782]"

$ws.Range("F11").Value = "[47 97] ; [15  5]"
$ws.Range("G11").Value = "[54.43902439 89.56097561] ; [ 7.56097561 12.43902439]"

$ws.Range("F12").Value = "[47 97] ; [15  5]"
$ws.Range("G12").Value = "[54.43902439 89.56097561] ; [ 7.56097561 12.43902439]"

$ws.Range("F13").Value = "[ 96 124] ; [17  3]"
$ws.Range("G13").Value = "[103.58333333 116.41666667] ; [ 9.41666667 10.58333333]"

$ws.Range("F14").Value = "[123 100] ; [19  1]"
$ws.Range("G14").Value = "[130.3127572  92.6872428] ; [11.6872428  8.3127572]"
